$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellAddr, $refAddr, $value) {
    $rng = $ws.Range($cellAddr)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = $ws.Range($refAddr).Style
}

Set-TextValue "D2" "B2" "274.15"
Set-TextValue "E2" "B2" "-1.77%"
Set-TextValue "D3" "B3" "26.65"
Set-TextValue "E3" "B3" "-2.69%"
Set-TextValue "D4" "B4" "4.763"
Set-TextValue "E4" "B4" "-0.72%"
Set-TextValue "D5" "B5" "0.06286"
Set-TextValue "E5" "B5" "-1.02%"
Set-TextValue "D6" "B6" "6.914"
Set-TextValue "E6" "B6" "-0.40%"
Set-TextValue "D7" "B7" "1.312"
Set-TextValue "E7" "B7" "37.26%"
Set-TextValue "D8" "B8" "0.8698"
Set-TextValue "E8" "B8" "-1.48%"
Set-TextValue "D9" "B9" "0.1573"
Set-TextValue "E9" "B9" "6.12%"
Set-TextValue "D10" "B10" "0.05027"
Set-TextValue "E10" "B10" "-4.47%"
Set-TextValue "D11" "B11" "0.07492"
Set-TextValue "E11" "B11" "2.79%"
Set-TextValue "D12" "B12" "0.02906"
Set-TextValue "E12" "B12" "-7.26%"
Set-TextValue "E13" "B13" "-0.09%"
Set-TextValue "D14" "B14" "0.001576"
Set-TextValue "E14" "B14" "0.82%"
Set-TextValue "D15" "B15" "0.0006352"
Set-TextValue "E15" "B15" "1.60%"
Set-TextValue "D16" "B16" "0.005838"
Set-TextValue "E16" "B16" "0.34%"
Set-TextValue "D17" "B17" "3.452"
Set-TextValue "E17" "B17" "-0.34%"
Set-TextValue "D18" "B18" "3.311"
Set-TextValue "E18" "B18" "-1.68%"
Set-TextValue "E19" "B19" "0.27%"
Set-TextValue "E20" "B20" "0.71%"
Set-TextValue "D21" "B21" "0.1317"
Set-TextValue "E21" "B21" "-1.62%"
Set-TextValue "D22" "B22" "3.933"
Set-TextValue "E22" "B22" "1.69%"
Set-TextValue "D23" "B23" "0.04391"
Set-TextValue "E23" "B23" "1.87%"
Set-TextValue "D24" "B24" "0.001169"
Set-TextValue "E24" "B24" "-0.90%"
Set-TextValue "D26" "B26" "0.0001201"
Set-TextValue "E26" "B26" "0.18%"
Set-TextValue "D27" "B27" "0.0001617"
Set-TextValue "E27" "B27" "-4.23%"
Set-TextValue "D40" "B40" "0.04074"
Set-TextValue "E40" "B40" "-0.43%"
Set-TextValue "D41" "B41" "0.007067"
Set-TextValue "E41" "B41" "5.65%"
Set-TextValue "D42" "B42" "0.1170"
Set-TextValue "E42" "B42" "0.59%"
Set-TextValue "D43" "B43" "0.002022"
Set-TextValue "E43" "B43" "-11.56%"
Set-TextValue "D44" "B44" "0.01122"
Set-TextValue "E44" "B44" "-10.27%"
Set-TextValue "D45" "B45" "0.00005199"
Set-TextValue "E45" "B45" "-0.44%"
Set-TextValue "D46" "B46" "0.02303"
Set-TextValue "E46" "B46" "2.40%"
Set-TextValue "D47" "B47" "1.490"
Set-TextValue "E47" "B47" "-37.38%"
